$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PasteSpecial constants used below.
$xlPasteAll = -4104
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# This mirrors the upstream content edit to docs/assets/disciplinas/LOT2028.xlsx:
# rows 13-24 get new text (several entries were re-ordered / replaced /
# removed), row heights are adjusted accordingly, and the former last row
# (25) is deleted outright so the sheet now ends at row 24.
#
# NOTE: a couple of the "new" cells end up (per the source diff) duplicating
# text that already lives elsewhere on the sheet (the activation date
# "01/01/2018" reappears in B15/C15, and the teacher name reappears in
# B18/C18). Those are captured here by copying the donor cell (value +
# format) first, which also avoids the text being mis-parsed back into a
# date value.
# ---------------------------------------------------------------------------

# Row 13: now carries the "Docentes responsáveis" teacher name in column A,
# plus a new "Semestral" entry in B/C.  Grab B18/C18's eventual text+style
# from the CURRENT B13/C13 before we overwrite them below. (The paste is
# done twice - once to get the value across safely, once more to stamp the
# correct cell format, since freshly-populated cells in a previously
# single-column row otherwise inherit column A's style.)
$ws.Range("B13").Copy()
$ws.Range("B18").PasteSpecial($xlPasteAll)
$ws.Range("B13").Copy()
$ws.Range("B18").PasteSpecial($xlPasteFormats)
$ws.Range("C13").Copy()
$ws.Range("C18").PasteSpecial($xlPasteAll)
$ws.Range("C13").Copy()
$ws.Range("C18").PasteSpecial($xlPasteFormats)

$ws.Rows.Item(13).RowHeight = 60
$ws.Range("A13").Value = "3403572 - Ismael Maciel de Mancilha"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 14: "Short syllabus:" row shifts up into row 14.
$ws.Rows.Item(14).RowHeight = 60
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Biotechnology (field of applications); fermentative processes; biochemistry of the fermentations (metabolic pathways of industrial interest); fermentative processes of industrial interest"
$ws.Range("C14").Value = "Biotechnology (field of applications); fermentative processes; biochemistry of the fermentations (metabolic pathways of industrial interest); fermentative processes of industrial interest"

# Row 15: "Programa:" label, with B/C now (per the source content) reusing
# the activation-date text. Copy straight from B8/C8 (value + format) so it
# stays text instead of being parsed into a date serial, then re-apply the
# formats to make sure the style matches column B/C's usual style exactly.
$ws.Rows.Item(15).RowHeight = 120
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial($xlPasteAll)
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial($xlPasteFormats)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial($xlPasteAll)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial($xlPasteFormats)

# Row 16: unchanged - "Syllabus:" plus the long English syllabus text.
$ws.Rows.Item(16).RowHeight = 120
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1.Biotechnology: concepts, application areas, multidisciplinary characteristic and examples of biotechnological products and processes.2.Fermentative processes: concept, enzymatic and fermentative processes, steps of fermentative process (downstream x upstream). Fermentative process modalities: a) batch and fed-batch fermentation, semi continuous and, continuous processes; b) induced and spontaneous fermentation;  c) semi solid fermentation;  d) oxygen supply; e) submerged and in surface processes; f) kinetics of the product formation in relation to the primary metabolism according to Gaden.3.Biochemistry of the fermentation: Fermentation – concepts, objectives, aerobic x anaerobic metabolisms; energy balance; preliminary steps of fermentation (extracellular hydrolysis and membrane permeability); metabolic pathways of industrial interest: a) EMP pathway; reactions and  allosteric control; alcoholic fermentation, homolactic fermentation, acetone/butanol, mixed-acid and 2,3 butanediol; b) Fosfo-Ketolase pathway; heterolactic fermentation and c) Entner Doudoroff pathway: alcoholic fermentation by Zymmonas mobilis. Fermentation balance: % recovered carbon and oxi-redox balance; Evaluation parameters of a fermentative process: yield, fermentation efficiency and productivity. Processes of interest: cocoa processing, ethanol production, fermented food and others."
$ws.Range("C16").Value = "1.Biotechnology: concepts, application areas, multidisciplinary characteristic and examples of biotechnological products and processes.2.Fermentative processes: concept, enzymatic and fermentative processes, steps of fermentative process (downstream x upstream). Fermentative process modalities: a) batch and fed-batch fermentation, semi continuous and, continuous processes; b) induced and spontaneous fermentation;  c) semi solid fermentation;  d) oxygen supply; e) submerged and in surface processes; f) kinetics of the product formation in relation to the primary metabolism according to Gaden.3.Biochemistry of the fermentation: Fermentation – concepts, objectives, aerobic x anaerobic metabolisms; energy balance; preliminary steps of fermentation (extracellular hydrolysis and membrane permeability); metabolic pathways of industrial interest: a) EMP pathway; reactions and  allosteric control; alcoholic fermentation, homolactic fermentation, acetone/butanol, mixed-acid and 2,3 butanediol; b) Fosfo-Ketolase pathway; heterolactic fermentation and c) Entner Doudoroff pathway: alcoholic fermentation by Zymmonas mobilis. Fermentation balance: % recovered carbon and oxi-redox balance; Evaluation parameters of a fermentative process: yield, fermentation efficiency and productivity. Processes of interest: cocoa processing, ethanol production, fermented food and others."

# Row 17: "Avaliação:" label only now - B/C are removed entirely, and the
# row height goes back to the sheet default (no customHeight flag).
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17:C17").Clear()
$ws.Rows.Item(17).AutoFit()

# Row 18: "Método:" label. B18/C18 were already populated above (copied
# from the original B13/C13 "3403572 - Ismael Maciel de Mancilha" cells),
# so just set the row height and column A text here.
$ws.Rows.Item(18).RowHeight = 60
$ws.Range("A18").Value = "Método:"

# Row 19: "Critério:" label with the evaluation-method text.
$ws.Rows.Item(19).RowHeight = 60
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "A avaliação será feita por meio de provas escritas."
$ws.Range("C19").Value = "A avaliação será feita por meio de provas escritas."

# Row 20: "Norma de recuperação:" label with the final-grade formula text.
$ws.Rows.Item(20).RowHeight = 60
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2) / 2"
$ws.Range("C20").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2) / 2"

# Row 21: "Bibliografia:" label with the recovery-rule text.
$ws.Rows.Item(21).RowHeight = 120
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"

# Row 22: "Requisitos:" label only - B/C removed entirely (the long
# bibliography text is dropped altogether), row height back to default.
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22:C22").Clear()
$ws.Rows.Item(22).AutoFit()

# Row 23: first prerequisite line moves up to row 23 (no column A anymore).
$ws.Range("A23").Clear()
$ws.Range("B23").Value = "LOT2008 -  Bioquímica II  (Requisito fraco)" + [char]10
$ws.Range("B10").Copy()
$ws.Range("B23").PasteSpecial($xlPasteFormats)
$ws.Range("C23").Value = "LOT2008 -  Bioquímica II  (Requisito fraco)" + [char]10
$ws.Range("C10").Copy()
$ws.Range("C23").PasteSpecial($xlPasteFormats)
$ws.Rows.Item(23).RowHeight = 30

# Row 24: second prerequisite line moves up to row 24.
$ws.Range("B24").Value = "LOT2053 -  Microbiologia  (Requisito fraco)" + [char]10
$ws.Range("C24").Value = "LOT2053 -  Microbiologia  (Requisito fraco)" + [char]10
$ws.Rows.Item(24).RowHeight = 30

# The old row 25 (second prerequisite line) is now gone - delete it so the
# sheet ends at row 24 and the dimension shrinks accordingly.
$ws.Rows.Item(25).Delete()
